# Ejercicio #2 complexity-analysis fix:
#   - The loop in se3_2 only runs to longitud/2, so its cost and the costs
#     of the instructions inside it are O(n/2) and O(n/2 - 1), not O(n) / O(n-1).
#   - The resulting Big-O derivation below the table is rewritten to match.

$d = $word.ActiveDocument

# --- Table 2 ("Ejercicio #2") - "Cuantas veces se repite" column ---
$t = $d.Tables.Item(2)
$t.Cell(3, 3).Range.Text = "O(n/2)"        # C2 row: for (...; i < longitud/2; ...)
$t.Cell(4, 3).Range.Text = "O(n/2 - 1)"    # C3 row: if (...)
$t.Cell(5, 3).Range.Text = "O(n/2 - 1)"    # C4 row: return "No es un palindromo";

# --- Formula paragraphs right after the table ---
# Touching a Table/Cell range leaves $d.Paragraphs' index stale in this
# runtime, so re-derive the paragraph collection from $d.Content instead of
# reusing $d.Paragraphs directly after the table edits above.
$body = $d.Content

# O(n) = C1 + C2 + C3(n-1) +C4(n-1) +C5  ->  ... C3(n/2 - 1) +C4(n/2 - 1) ...
$body.Paragraphs.Item(97).Range.Text = "O(n) = C1 + C2 + C3(n/2 - 1) +C4(n/2 - 1) +C5"

# Drop the intermediate expansion step entirely:
#   "O(n) = C1 + C2 +C3n - C3 + C4n - C4 + C5 "
$d.Content.Paragraphs.Item(98).Range.Delete()

# The final simplified-formula paragraph (now shifted up to index 98) is
# rewritten with the corrected, expanded closed form.
$d.Content.Paragraphs.Item(98).Range.Text = "O(n) = (C1 + C5 - C3 - C4) + n(C2 + C3 +C4) + ½(C2 + C3 + C4)"
